$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" conversion note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.69 = 9989.5 pesos`n✅ 9989.5 pesos = 2.69 = 966.98 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" exchange-rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 371.59
$wsTasas.Range("O10").Value = 3712
$wsTasas.Range("N12").Value = 3720
$wsTasas.Range("O12").Value = 360.096
